$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds the same serial date value (45206) for every
# data row (rows 2-238). The edit bumps that date from 45206 to 45208
# (2023-10-07 -> 2023-10-09) for all of them.
$ws.Range("C2:C238").Value = 45208
